$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append the new inventory items to the bottom of the single-column list
# (sheet1, column A). Previously the sheet ran from A1 (header "Item")
# through A308 ("Incense Stick (Agarbatthi)"). We're adding six more rows.
# ---------------------------------------------------------------------------

$newItems = @(
    "Coffee_brass_set",
    "Ulavacharu",
    "Milk_bread",
    "pen",
    "battery",
    "silver_pouch"
)

$firstNewRow = 309
$row = $firstNewRow
foreach ($item in $newItems) {
    $ws.Cells.Item($row, 1).Value = $item
    $row++
}
$lastNewRow = $row - 1

# ---------------------------------------------------------------------------
# Formatting: the existing tail of the list (rows 296:308) carried a
# "last block" cell style (thin left/right border, no fill). Adding new
# rows below pushes that visual treatment down onto the newly added rows,
# while the rows that used to be last fall back to the regular list style.
# Both styles render identically (thin left/right border, no fill) - we
# simply reapply the border formatting so the style bookkeeping matches.
# ---------------------------------------------------------------------------

$oldTailRange = $ws.Range("A296:A308")
$oldTailRange.Borders.Item(7).LineStyle = 1    # xlEdgeLeft
$oldTailRange.Borders.Item(7).Weight = 2       # xlThin
$oldTailRange.Borders.Item(10).LineStyle = 1   # xlEdgeRight
$oldTailRange.Borders.Item(10).Weight = 2      # xlThin

$newRange = $ws.Range("A" + $firstNewRow + ":A" + $lastNewRow)
$newRange.Borders.Item(7).LineStyle = 1
$newRange.Borders.Item(7).Weight = 2
$newRange.Borders.Item(10).LineStyle = 1
$newRange.Borders.Item(10).Weight = 2

# ---------------------------------------------------------------------------
# Update the active selection to point at the new last cell, matching
# where Excel would leave the cursor after typing the new entries.
# ---------------------------------------------------------------------------

$nextCell = $ws.Range("A" + ($lastNewRow + 1))
$nextCell.Select() | Out-Null
